$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.454385280609131
$ws.Range("B1").Value = 2.980066061019897
$ws.Range("C1").Value = 2.622552633285522
$ws.Range("D1").Value = 1.670663833618164
$ws.Range("E1").Value = 0.8170940279960632
